$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 312, shifting existing rows 312:399 down to 313:400
$ws.Rows("312:312").Insert()

# Populate the newly inserted row 312 with the new data record
$ws.Range("A312").Value = 10
$ws.Range("B312").Value = "Vega Modelo de Temuco"
$ws.Range("C312").Value = "La Araucanía"
$ws.Range("D312").Value = 44663
$ws.Range("E312").Value = 9
$ws.Range("F312").Value = "Fruta"
$ws.Range("G312").Value = 100101
$ws.Range("H312").Value = "Berries"
$ws.Range("I312").Value = 100101007
$ws.Range("J312").Value = "Kiwi"
$ws.Range("K312").Value = "Hayward"
$ws.Range("L312").Value = "Primera"
$ws.Range("M312").Value = 155
$ws.Range("N312").Value = 12000
$ws.Range("O312").Value = 12000
$ws.Range("P312").Value = 12000
$ws.Range("Q312").Value = "$/bandeja 10 kilos"
$ws.Range("R312").Value = "Región de O'Higgins"
$ws.Range("S312").Value = 1200
$ws.Range("T312").Value = 10
